$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "folder_name" column (column D). This shifts columns E..J left to D..I.
$ws.Columns.Item(4).Delete()

# Add the new third row of data.
$ws.Range("A3").Value = "DemoCollection2.postman_collection.json"
$ws.Range("B3").Value = "createuser/Demodata.json"
$ws.Range("E3").Value = "createuser/QA2.postman_environment.json"
$ws.Range("G3").Value = "RecieveDataRequests"
$ws.Range("H3").Value = "ReceiveDataRequests.html"
$ws.Range("I3").Value = "receivedata"

# Update the active selection to match the target state.
$ws.Range("B3").Select()
